# Actualización 10 de Mayo
# Updates the Aprobados/Reprobados/Por_Apro/Por_Repro/Promedio/Blancos/Por_Blan
# columns (E-K) for the first 5 student groups on each of the three exam sheets.

$wb = $excel.ActiveWorkbook

# --- Hoja "1er Parcial" ---
$ws = $wb.Worksheets.Item("1er Parcial")
$ws.Range("I2").Value = 6.3
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

$ws.Range("E3").Value = 23
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = 74.19
$ws.Range("H3").Value = 25.81
$ws.Range("I3").Value = 6.4
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0

$ws.Range("E5").Value = 34
$ws.Range("F5").Value = 6
$ws.Range("G5").Value = 85
$ws.Range("H5").Value = 15
$ws.Range("I5").Value = 7.6
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0

$ws.Range("I6").Value = 6.4
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0

# --- Hoja "2o Parcial" ---
$ws = $wb.Worksheets.Item("2o Parcial")
$ws.Range("E2").Value = 14
$ws.Range("F2").Value = 22
$ws.Range("G2").Value = 38.89
$ws.Range("H2").Value = 61.11
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

$ws.Range("E3").Value = 21
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 67.74
$ws.Range("H3").Value = 32.26
$ws.Range("I3").Value = 7.3
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

$ws.Range("E4").Value = 14
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 66.67
$ws.Range("H4").Value = 33.33
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0

$ws.Range("E5").Value = 30
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = 75
$ws.Range("H5").Value = 25
$ws.Range("I5").Value = 7.3
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0

$ws.Range("E6").Value = 13
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 56.52
$ws.Range("H6").Value = 43.48
$ws.Range("I6").Value = 5.8
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0

# --- Hoja "3er Parcial" ---
$ws = $wb.Worksheets.Item("3er Parcial")
$ws.Range("E2").Value = 14
$ws.Range("F2").Value = 22
$ws.Range("G2").Value = 38.89
$ws.Range("H2").Value = 61.11
$ws.Range("I2").Value = 6.1
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

$ws.Range("E3").Value = 21
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 67.74
$ws.Range("H3").Value = 32.26
$ws.Range("I3").Value = 6.9
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

$ws.Range("E4").Value = 14
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 66.67
$ws.Range("H4").Value = 33.33
$ws.Range("I4").Value = 7.1
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0

$ws.Range("E5").Value = 30
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = 75
$ws.Range("H5").Value = 25
$ws.Range("I5").Value = 7.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0

$ws.Range("E6").Value = 13
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 56.52
$ws.Range("H6").Value = 43.48
$ws.Range("I6").Value = 6.2
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
